$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.434.43'
$ws.Range("E2").Value = '  -1.65%  '

$ws.Range("D3").Value = '1.795.95'
$ws.Range("E3").Value = '  -1.93%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("E5").Value = '  +0.11%  '

$ws.Range("D6").Value = "'307.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.92%  '

$ws.Range("D7").Value = "'0.4539"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.57%  '

$ws.Range("D8").Value = "'0.3596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.03%  '

$ws.Range("D9").Value = "'46.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.12%  '

$ws.Range("D10").Value = "'0.07124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.60%  '

$ws.Range("D11").Value = "'0.8884"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.34%  '

$ws.Range("D12").Value = "'0.07823"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.87%  '

$ws.Range("D13").Value = "'19.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.26%  '

$ws.Range("D14").Value = '1.765.52'
$ws.Range("E14").Value = '  -3.80%  '

$ws.Range("D15").Value = "'5.283"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.99%  '

$ws.Range("D16").Value = "'6.339"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.65%  '

$ws.Range("E17").Value = '  -2.76%  '

$ws.Range("D18").Value = "'1.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.25%  '

$ws.Range("D19").Value = "'0.000008589"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.54%  '

$ws.Range("E20").Value = '  +0.14%  '

$ws.Range("D21").Value = "'14.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.01%  '

$ws.Range("D22").Value = '26.446.69'
$ws.Range("E22").Value = '  -1.76%  '

$ws.Range("D23").Value = "'4.992"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("E24").Value = '  +1.22%  '

$ws.Range("D25").Value = '2.001.11'
$ws.Range("E25").Value = '  -4.26%  '

$ws.Range("D26").Value = "'1.986"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").Value = "'152.71"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").Value = "'17.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.56%  '

$ws.Range("D29").Value = "'2.049"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.21%  '

$ws.Range("D30").Value = "'112.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.19%  '

$ws.Range("D31").Value = "'4.877"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.12%  '

$ws.Range("D32").Value = "'0.08661"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.99%  '

$ws.Range("D33").Value = "'3.053"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.38%  '

$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = "'2.757"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.79%  '

$ws.Range("D35").Value = "'0.7304"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.08%  '

$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").Value = "'4.457"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("E37").Value = '  -1.42%  '

$ws.Range("E38").Value = '  -1.00%  '

$ws.Range("D39").Value = "'0.01943"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.36%  '

$ws.Range("D40").Value = "'0.05118"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.17%  '

$ws.Range("D41").Value = "'2.876"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.74%  '

$ws.Range("D42").Value = "'0.5148"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.44%  '

$ws.Range("D43").Value = "'6.894"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '

$ws.Range("E44").Value = '  -4.71%  '

$ws.Range("D45").Value = "'8.007"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.83%  '

$ws.Range("E46").Value = '  +0.14%  '

$ws.Range("D47").Value = "'0.4669"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.10%  '

$ws.Range("D48").Value = "'9.909"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("D49").Value = "'100.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.62%  '

$ws.Range("D50").Value = "'1.588"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.36%  '

$ws.Range("D51").Value = "'0.05984"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.83%  '
